$wb = $excel.ActiveWorkbook

# Sheet "展览"
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 10596
$wsExhibit.Range("F5").Value = 670
$wsExhibit.Range("F6").Value = 490

# Sheet "全部类型"
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 10596
$wsAll.Range("F5").Value = 670
$wsAll.Range("F7").Value = 490
